# Build site at 2021-10-27 12:19:33 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Semestre ideal:" value EQD-9,EQN-12 -> EQD-9,EQN-11 (both display columns)
$ws.Range("B9").Value = "EQD-9,EQN-11"
$ws.Range("C9").Value = "EQD-9,EQN-11"

# Remove two of the four "Requisitos" rows:
#   row 24 = LOM3081 - Introdução à Mecânica dos Sólidos (Requisito fraco)
#   row 26 = LOQ4054 - Fenômenos de Transporte III (Requisito fraco)
# leaving LOQ4002 and LOQ4086, which shift up to rows 24 and 25.
$ws.Rows("24").Delete()
$ws.Rows("25").Delete()
